$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting from the last existing data row (29) down into the
# three new rows so the new cells reuse the same cell style (date format
# on column A, etc.) instead of creating new style entries.
$ws.Range("A29:M29").Copy() | Out-Null
$ws.Range("A30:M32").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$data = @(
    @(42602.576064814813, "Bag", 163, 338, 31, 7, 16, 30, 69, 0, 3, 0, 99),
    @(42602.576898148145, "Bag", 144, 338, 31, 7, 16, 30, 69, 0, 3, 0, 99),
    @(42602.577210648145, "Bag", 138, 338, 31, 4, 18, 18, 81, 0, 4, 0, 100)
)

$startRow = 30
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 1; $c -le 13; $c++) {
        $ws.Cells.Item($row, $c).Value = $values[$c - 1]
    }
}
